# Workbook / worksheet handles
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rename the sheet (Sheet1 -> samplefile) ---
$ws.Name = "samplefile"

# --- Insert a new column before column A; this shifts the existing
#     columns A:J (headers + data) one position to the right, to B:K ---
$ws.Columns.Item(1).Insert()

# --- New "date" header in A1 ---
$ws.Range("A1").Value = "date"

# --- Apply the built-in date number format (numFmtId 15, "d-mmm-yy")
#     to the new date column before filling it in ---
$ws.Range("A2:A21").NumberFormat = "d-mmm-yy"

# --- Fill A2:A21 with the sequential dates 2019-05-21 .. 2019-06-09 ---
$baseDate = Get-Date -Year 2019 -Month 5 -Day 21 -Hour 0 -Minute 0 -Second 0
for ($i = 0; $i -lt 20; $i++) {
    $ws.Cells.Item(2 + $i, 1).Value = $baseDate.AddDays($i)
}

# --- Autofit the new column to its content ---
$ws.Columns.Item(1).AutoFit()

# --- Restore the selection highlighted in the sheet ---
$ws.Range("I9:I14").Select()
